# Add a new "Slovakia" worksheet (based on the "Portugal" sheet's layout)
# to the workbook, fill it in with the Slovakia market data, and update
# the sheet selections so that the new sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Portugal is no longer the tab that is shown/selected - simulate the
# user clicking on the row 1 header (selecting the whole sheet) before
# switching away from it.
$portugal.Cells.Select()

# Duplicate the Portugal sheet (keeps formatting, column widths, merged
# cells, page setup, etc.) and place the copy at the end of the tab strip.
$portugal.Copy($null, $portugal)

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Update the market specific values.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222"

# The copied rows inherited Portugal's explicit (wrap-text driven) row
# heights; re-fit them back down to the sheet's natural/default height,
# matching a freshly authored sheet.
$slovakia.Rows("1:15").AutoFit()

# Make the new sheet the active one, with the same cell selected as in
# the authored workbook.
$slovakia.Activate()
$slovakia.Range("C17").Select()
